$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content corrections (Use Case Scenario: Kreiranje korisničkog računa) ---

# "Preduvjeti" (Preconditions): the user no longer needs Internet access,
# instead they must go to the Gras counter to give data to the administrator.
$ws.Range("B3").Value = "Korisnik mora doći na Grasov šalter kako bi dao svoje podatke administratoru"

# "Posljedice - neuspjesan zavrsetak" (Consequences - unsuccessful outcome):
# dropped the re-entry clause.
$ws.Range("B5").Value = "Korisnik dobiva obavijest da njegov račun nije kreiran "

# "Glavni tok" (Main flow) summary now mentions the administrator doing the work.
$ws.Range("B7").Value = "Korisnik podnosi zahtjev za formiranje korisničkog računa, administrator popuni odgovarajuće podatke, otvori korisnički račun i dodjeli korisničko ime i password"

# "Prosirenja/Alternative" (Extensions/Alternatives): now the administrator
# (not the user) can mistype the data.
$ws.Range("B8").Value = "Administrator pogrešno unese podatke, sistem traži ponovni unos podataka"

# Steps of the main flow of events: step 3, 5 and 6 are replaced with new,
# plain (non rich-text) descriptions reflecting the administrator's role.
$ws.Range("A18").Value = "3. Odlazak na Grasov šalter"
$ws.Range("B20").Value = "5. Upis korisnikovih podataka"
$ws.Range("B21").Value = "6. Dodjela korisničkog imena i password-a"

# --- Row height adjustments to fit the revised text ---
$ws.Rows.Item(3).RowHeight = 29.25
$ws.Rows.Item(7).RowHeight = 60.75

# --- View state: scroll down and select B21:B22 ---
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B21:B22").Select() | Out-Null
